$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph that currently sits
# right after the title heading (it gets moved/reworded down near the end
# of the document instead). ---
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: the old last paragraph held an (unused) AI image-prompt.
# Replace it with two new paragraphs: a bold title line, followed by the
# meta-description text (now italic, without the "Meta description: "
# label prefix). ---
$targetPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Please create a feature image*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    $targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
}

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Gangsterz Free Slot Game | Review 2021</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Gangsterz Slot by Top Trend Gaming. Find out how to play and win in this cluster-based game. Play Gangsterz free and enjoy high-quality graphics!</w:t></w:r></w:p>'
$targetPara.Range.InsertXML($newXml)
